$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16357.368
$ws.Range("I69").Value = 14166.667
$ws.Range("J69").Value = 16768.125
$ws.Range("K69").Value = 42500.001
$ws.Range("L69").Value = 50304.375
$ws.Range("M69").Value = -41626.001
$ws.Range("N69").Value = -52052.375

$ws.Range("H72").Value = 16357.368
$ws.Range("I72").Value = 14166.667
$ws.Range("J72").Value = 16768.125
$ws.Range("K72").Value = 127500.003
$ws.Range("L72").Value = 150913.125
$ws.Range("M72").Value = -123132.003
$ws.Range("N72").Value = -159649.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1429.2667
$ws.Range("I110").Value = 1367.619
$ws.Range("J110").Value = 1573.1111
$ws.Range("K110").Value = 1367.619
$ws.Range("L110").Value = 1573.1111
$ws.Range("M110").Value = 677.3810000000001
$ws.Range("N110").Value = -5663.1111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1541.5652
$ws.Range("I20").Value = 1352.8334
$ws.Range("J20").Value = 1747.4546
$ws.Range("K20").Value = 1352.8334
$ws.Range("L20").Value = 1747.4546
$ws.Range("M20").Value = -1105.8334
$ws.Range("N20").Value = -2241.4546

$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4259.6
$ws.Range("I15").Value = 2800
$ws.Range("J15").Value = 5232.6665
$ws.Range("K15").Value = 2800
$ws.Range("L15").Value = 5232.6665
$ws.Range("M15").Value = -2630
$ws.Range("N15").Value = -5572.6665

$ws.Range("H31").Value = 1645.0212
$ws.Range("I31").Value = 992.5599999999999
$ws.Range("J31").Value = 2386.4546
$ws.Range("K31").Value = 992.5599999999999
$ws.Range("L31").Value = 2386.4546
$ws.Range("M31").Value = -697.5599999999999
$ws.Range("N31").Value = -2976.4546

$ws.Range("H34").Value = 1645.0212
$ws.Range("I34").Value = 992.5599999999999
$ws.Range("J34").Value = 2386.4546
$ws.Range("K34").Value = 992.5599999999999
$ws.Range("L34").Value = 2386.4546
$ws.Range("M34").Value = -790.5599999999999
$ws.Range("N34").Value = -2790.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1200.3334
$ws.Range("I68").Value = 1150.6666
$ws.Range("J68").Value = 1250
$ws.Range("K68").Value = 3451.9998
$ws.Range("L68").Value = 3750
$ws.Range("M68").Value = -2640.9998
$ws.Range("N68").Value = -5372

$ws.Range("H71").Value = 1200.3334
$ws.Range("I71").Value = 1150.6666
$ws.Range("J71").Value = 1250
$ws.Range("K71").Value = 10355.9994
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -6299.999400000001
$ws.Range("N71").Value = -19362

$ws.Range("H129").Value = 3738.4736
$ws.Range("I129").Value = 5277.5
$ws.Range("J129").Value = 3557.4119
$ws.Range("K129").Value = 15832.5
$ws.Range("L129").Value = 10672.2357
$ws.Range("M129").Value = -10832.5
$ws.Range("N129").Value = -20672.2357

$ws.Range("H131").Value = 899.65
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 920.6842
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2762.0526
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12842.0526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1573055.9
$ws.Range("I3").Value = 1833565.1
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 1833565.1
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -1833449.1
$ws.Range("N3").Value = -10232

$ws.Range("H70").Value = 5510.5625
$ws.Range("I70").Value = 5116.364
$ws.Range("J70").Value = 6377.8
$ws.Range("K70").Value = 5116.364
$ws.Range("L70").Value = 6377.8
$ws.Range("M70").Value = -4846.364
$ws.Range("N70").Value = -6917.8

$ws.Range("H73").Value = 5510.5625
$ws.Range("I73").Value = 5116.364
$ws.Range("J73").Value = 6377.8
$ws.Range("K73").Value = 5116.364
$ws.Range("L73").Value = 6377.8
$ws.Range("M73").Value = -4180.364
$ws.Range("N73").Value = -8249.799999999999

$ws.Range("H122").Value = 173392.5
$ws.Range("I122").Value = 253700
$ws.Range("J122").Value = 12777.5
$ws.Range("K122").Value = 761100
$ws.Range("L122").Value = 38332.5
$ws.Range("M122").Value = -758650
$ws.Range("N122").Value = -43232.5

$ws.Range("H132").Value = 2566883
$ws.Range("I132").Value = 2723.2068
$ws.Range("J132").Value = 10002946
$ws.Range("K132").Value = 8169.6204
$ws.Range("L132").Value = 30008838
$ws.Range("M132").Value = -5639.6204
$ws.Range("N132").Value = -30013898

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2066.3333
$ws.Range("I40").Value = 2042.7142
$ws.Range("J40").Value = 2113.5715
$ws.Range("K40").Value = 2042.7142
$ws.Range("L40").Value = 2113.5715
$ws.Range("M40").Value = -1906.7142
$ws.Range("N40").Value = -2385.5715

$ws.Range("H46").Value = 913.9583
$ws.Range("I46").Value = 1108.25
$ws.Range("J46").Value = 816.8125
$ws.Range("K46").Value = 1108.25
$ws.Range("L46").Value = 816.8125
$ws.Range("M46").Value = -920.25
$ws.Range("N46").Value = -1192.8125

$ws.Range("H68").Value = 10646.077
$ws.Range("I68").Value = 18750.666
$ws.Range("J68").Value = 3699.2856
$ws.Range("K68").Value = 18750.666
$ws.Range("L68").Value = 3699.2856
$ws.Range("M68").Value = -18001.666
$ws.Range("N68").Value = -5197.2856

$ws.Range("H71").Value = 10646.077
$ws.Range("I71").Value = 18750.666
$ws.Range("J71").Value = 3699.2856
$ws.Range("K71").Value = 93753.33
$ws.Range("L71").Value = 18496.428
$ws.Range("M71").Value = -90009.33
$ws.Range("N71").Value = -25984.428

$ws.Range("H82").Value = 1444.04
$ws.Range("I82").Value = 1681
$ws.Range("J82").Value = 1225.3077
$ws.Range("K82").Value = 1681
$ws.Range("L82").Value = 1225.3077
$ws.Range("M82").Value = -1320
$ws.Range("N82").Value = -1947.3077

$ws.Range("H85").Value = 1444.04
$ws.Range("I85").Value = 1681
$ws.Range("J85").Value = 1225.3077
$ws.Range("K85").Value = 1681
$ws.Range("L85").Value = 1225.3077
$ws.Range("M85").Value = -433
$ws.Range("N85").Value = -3721.3077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H62").Value = 2950
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2326
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2950
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11630
$ws.Range("N65").ClearContents()

$ws.Range("H81").Value = 1240.4286
$ws.Range("I81").Value = 1240.4286
$ws.Range("K81").Value = 2480.8572
$ws.Range("M81").Value = -1419.8572

$ws.Range("H84").Value = 1240.4286
$ws.Range("I84").Value = 1240.4286
$ws.Range("K84").Value = 12404.286
$ws.Range("M84").Value = -7100.286
